$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.260.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.676.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5275'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.65%  '
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2659'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06297'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07569'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.673.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.477'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5642'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.31%  '
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008049'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.084.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.831'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '188.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.215'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.004'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '150.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1258'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.12%  '
$ws.Range('E26').Value = '  -3.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06228'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.363'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.287'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.514'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.449'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.638'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.33%  '
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6075'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.109.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01625'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.113'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8721'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.827.36'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000107'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.006'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.011'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05236'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.985'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.85%  '
